# Populate the expense-details sheet: replace the single sample "food" row
# with the real expense rows (Cafe/movie/travel/clothing/rent), extending
# the sheet from A1:C2 to A1:C6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update existing row
$ws.Range("A2").Value = "Cafe"
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = "22/5/2025"

# Row 3 - new row
$ws.Range("A3").Value = "movie"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = "16/5/2025"

# Row 4 - new row
$ws.Range("A4").Value = "travel"
$ws.Range("B4").Value = 3000
$ws.Range("C4").Value = "15/5/2025"

# Row 5 - new row
$ws.Range("A5").Value = "clothing"
$ws.Range("B5").Value = 16000
$ws.Range("C5").Value = "14/5/2025"

# Row 6 - new row
$ws.Range("A6").Value = "rent"
$ws.Range("B6").Value = 10000
# "10/5/2025" parses as a valid M/D/Y date (day=10 <= 12), so Excel would
# silently convert it to a date serial on assignment. Force text interpretation
# (like Excel's "Text" format / leading apostrophe would), then restore the
# cell to the default (unstyled) state so no stray style sticks around.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "10/5/2025"
$ws.Range("C6").Style = "Normal"
